# FU_TEMPLATE_Klickovich.docx update
#
# The "Insurance" line used two separate merge fields
#     Insurance:  {{insurance1}}  {{insurance2}}
# and is consolidated into a single list-producing field
#     Insurance:  {{insuranceList}}
#
# Doing that edit interactively in Word leaves the cursor right after
# "Insurance:  " / right before the new "{{" - which is exactly where
# real Word drops its "last edit location" bookmark (the hidden
# "_GoBack" bookmark). Recreating that bookmark here reproduces the
# bookmark-id bump that the diff shows on the pre-existing
# "_Hlk129781177" bookmark (0 -> 1), since "_GoBack" takes over id 0.

$d = $word.ActiveDocument

# 1) Collapse the two insurance fields into the single insuranceList field.
$d.Content.Find.Execute("{{insurance1}}  {{insurance2}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{insuranceList}}", 2) | Out-Null

# 2) Re-find the freshly inserted field so we can drop a collapsed
#    range right in front of it (immediately after "Insurance:  ").
$hit = $d.Content
$hit.Find.Execute("{{insuranceList}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$caret = $d.Range($hit.Start, $hit.Start)
$d.Bookmarks.Add("_GoBack", $caret) | Out-Null
